# Update the 江西-漫展信息 workbook: remove six outdated rows and refresh
# the "want to go" counts on the surviving rows, for both the "展览" and
# "全部类型" worksheets (which carry identical data).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Row numbers (in the *original* sheet) to delete. Delete from the bottom
# up so earlier deletions don't shift the row numbers we still need to
# remove.
$rowsToDelete = @(11, 10, 8, 6, 4, 3)

# F-column ("want to go" count) corrections to apply to the rows that
# remain *after* the deletions above have shifted everything up.
$fUpdates = @{
    3  = 8541
    4  = 1520
    6  = 391
    7  = 253
    8  = 169
    11 = 44
    12 = 462
    13 = 1245
    14 = 225
    16 = 138
    17 = 97
    18 = 127
    19 = 76
    20 = 114
    21 = 105
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($r in $rowsToDelete) {
        $ws.Rows.Item($r).Delete()
    }

    foreach ($r in $fUpdates.Keys) {
        $ws.Cells.Item($r, 6).Value = $fUpdates[$r]
    }
}
